# rerun corona results with larger ds
# ------------------------------------------------------------------
# This script mutates the already-open workbook ($excel.ActiveWorkbook)
# so that results/.../avg_0.004_scores.xlsx reflects a rerun of the
# scoring pass against a larger dataset: most counts/percentages in the
# two result tables (A1:H6 "negative" table, J1:Q30 "positive" table)
# change, the old word rows 7-9 ("sc"/"no"/"stop") disappear from the
# negative table, and four new ranked rows (27-30: help/increase/
# please/".") are appended to the bottom of the positive table.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The larger dataset only produced 6 ranked words for the "negative"
# table instead of 9 - fully remove the old A7:H9 block (not just its
# contents) so no stray styled/empty cells remain.
$ws.Range("A7:H9").Clear()

# Rows 3-26: refreshed counts/percentages + re-ranked words for both
# tables.
$ws.Range("B3").Value = 0.9117647058823529
$ws.Range("C3").Value = 31
$ws.Range("D3").Value = 31
$ws.Range("H3").Value = 3
$ws.Range("J3").Value = 'best'
$ws.Range("K3").Value = 0.9152542372881356
$ws.Range("L3").Value = 54
$ws.Range("M3").Value = 54
$ws.Range("Q3").Value = 5
$ws.Range("B4").Value = 0.5856164383561644
$ws.Range("C4").Value = 171
$ws.Range("D4").Value = 171
$ws.Range("H4").Value = 121
$ws.Range("J4").Value = 'interesting'
$ws.Range("K4").Value = 0.9090909090909091
$ws.Range("L4").Value = 30
$ws.Range("M4").Value = 30
$ws.Range("Q4").Value = 3
$ws.Range("B5").Value = 0.1744186046511628
$ws.Range("C5").Value = 90
$ws.Range("D5").Value = 90
$ws.Range("H5").Value = 426
$ws.Range("K5").Value = 0.8913043478260869
$ws.Range("L5").Value = 41
$ws.Range("M5").Value = 41
$ws.Range("Q5").Value = 5
$ws.Range("A6").Value = 'sc'
$ws.Range("B6").Value = 0.1693121693121693
$ws.Range("C6").Value = 32
$ws.Range("D6").Value = 32
$ws.Range("H6").Value = 157
$ws.Range("J6").Value = 'great'
$ws.Range("K6").Value = 0.8571428571428571
$ws.Range("L6").Value = 96
$ws.Range("M6").Value = 96
$ws.Range("Q6").Value = 16
$ws.Range("J7").Value = 'won'
$ws.Range("K7").Value = 0.8205128205128205
$ws.Range("L7").Value = 32
$ws.Range("M7").Value = 32
$ws.Range("J8").Value = 'positive'
$ws.Range("K8").Value = 0.7931034482758621
$ws.Range("L8").Value = 46
$ws.Range("M8").Value = 46
$ws.Range("Q8").Value = 12
$ws.Range("J9").Value = 'thank'
$ws.Range("K9").Value = 0.7890625
$ws.Range("L9").Value = 101
$ws.Range("M9").Value = 101
$ws.Range("Q9").Value = 27
$ws.Range("J10").Value = 'thanks'
$ws.Range("K10").Value = 0.7804878048780488
$ws.Range("L10").Value = 64
$ws.Range("M10").Value = 64
$ws.Range("Q10").Value = 18
$ws.Range("J11").Value = 'safe'
$ws.Range("K11").Value = 0.7394366197183099
$ws.Range("L11").Value = 105
$ws.Range("M11").Value = 105
$ws.Range("Q11").Value = 37
$ws.Range("J12").Value = 'support'
$ws.Range("K12").Value = 0.7358490566037735
$ws.Range("L12").Value = 78
$ws.Range("M12").Value = 78
$ws.Range("Q12").Value = 28
$ws.Range("J13").Value = 'free'
$ws.Range("K13").Value = 0.7333333333333333
$ws.Range("L13").Value = 88
$ws.Range("M13").Value = 88
$ws.Range("Q13").Value = 32
$ws.Range("J14").Value = 'special'
$ws.Range("K14").Value = 0.7222222222222222
$ws.Range("L14").Value = 26
$ws.Range("M14").Value = 26
$ws.Range("Q14").Value = 10
$ws.Range("J15").Value = 'good'
$ws.Range("K15").Value = 0.71875
$ws.Range("L15").Value = 115
$ws.Range("M15").Value = 115
$ws.Range("Q15").Value = 45
$ws.Range("J16").Value = 'safety'
$ws.Range("K16").Value = 0.7058823529411765
$ws.Range("L16").Value = 36
$ws.Range("M16").Value = 36
$ws.Range("Q16").Value = 15
$ws.Range("J17").Value = 'confidence'
$ws.Range("K17").Value = 0.6944444444444444
$ws.Range("L17").Value = 25
$ws.Range("M17").Value = 25
$ws.Range("Q17").Value = 11
$ws.Range("J18").Value = 'better'
$ws.Range("K18").Value = 0.6507936507936508
$ws.Range("L18").Value = 41
$ws.Range("M18").Value = 41
$ws.Range("Q18").Value = 22
$ws.Range("J19").Value = 'relief'
$ws.Range("K19").Value = 0.64
$ws.Range("L19").Value = 32
$ws.Range("M19").Value = 32
$ws.Range("Q19").Value = 18
$ws.Range("J20").Value = 'heroes'
$ws.Range("K20").Value = 0.6382978723404256
$ws.Range("L20").Value = 30
$ws.Range("M20").Value = 30
$ws.Range("Q20").Value = 17
$ws.Range("J21").Value = 'well'
$ws.Range("K21").Value = 0.6170212765957447
$ws.Range("L21").Value = 58
$ws.Range("M21").Value = 58
$ws.Range("Q21").Value = 36
$ws.Range("J22").Value = 'fresh'
$ws.Range("K22").Value = 0.5833333333333334
$ws.Range("L22").Value = 28
$ws.Range("M22").Value = 28
$ws.Range("Q22").Value = 20
$ws.Range("J23").Value = 'hand'
$ws.Range("K23").Value = 0.5169712793733682
$ws.Range("L23").Value = 198
$ws.Range("M23").Value = 198
$ws.Range("Q23").Value = 185
$ws.Range("J24").Value = 'care'
$ws.Range("K24").Value = 0.4606741573033708
$ws.Range("L24").Value = 41
$ws.Range("M24").Value = 41
$ws.Range("Q24").Value = 48
$ws.Range("J25").Value = 'like'
$ws.Range("K25").Value = 0.4441176470588235
$ws.Range("L25").Value = 151
$ws.Range("M25").Value = 151
$ws.Range("Q25").Value = 189
$ws.Range("J26").Value = 'hope'
$ws.Range("K26").Value = 0.4307692307692308
$ws.Range("L26").Value = 28
$ws.Range("M26").Value = 28
$ws.Range("Q26").Value = 37

# Rows 27-30 are brand new - copy the bold/bordered/centered format
# used by the existing word cells in column J (e.g. J26) before
# filling in the new ranked words + stats so the new cells match the
# table styling exactly.
$ws.Range("J26").Copy() | Out-Null
$ws.Range("J27").PasteSpecial(-4122) | Out-Null
$ws.Range("J27").Value = 'help'
$ws.Range("K27").Value = 0.4271186440677966
$ws.Range("L27").Value = 126
$ws.Range("M27").Value = 126
$ws.Range("N27").Value = 1
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = $false
$ws.Range("Q27").Value = 169
$ws.Range("J26").Copy() | Out-Null
$ws.Range("J28").PasteSpecial(-4122) | Out-Null
$ws.Range("J28").Value = 'increase'
$ws.Range("K28").Value = 0.3333333333333333
$ws.Range("L28").Value = 26
$ws.Range("M28").Value = 26
$ws.Range("N28").Value = 1
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = $false
$ws.Range("Q28").Value = 52
$ws.Range("J26").Copy() | Out-Null
$ws.Range("J29").PasteSpecial(-4122) | Out-Null
$ws.Range("J29").Value = 'please'
$ws.Range("K29").Value = 0.3221757322175732
$ws.Range("L29").Value = 77
$ws.Range("M29").Value = 77
$ws.Range("N29").Value = 1
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = $false
$ws.Range("Q29").Value = 162
$ws.Range("J26").Copy() | Out-Null
$ws.Range("J30").PasteSpecial(-4122) | Out-Null
$ws.Range("J30").Value = '.'
$ws.Range("K30").Value = 0.005804643714971977
$ws.Range("L30").Value = 29
$ws.Range("M30").Value = 29
$ws.Range("N30").Value = 1
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = $false
$ws.Range("Q30").Value = 4967

$excel.CutCopyMode = $false
